# Refresh cryptos list data (prices + 1h volume change) pulled from coinranking.com.
# Source rows 2-51 on Sheet1; columns: B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '99.018.77'
$ws.Range("E2").Value = '  +2.16%  '

$ws.Range("D3").Value = '3.401.98'
$ws.Range("E3").Value = '  +8.90%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''261.02'
$ws.Range("E5").Value = '  +8.88%  '

$ws.Range("D6").Value = '''636.39'
$ws.Range("E6").Value = '  +4.64%  '

$ws.Range("E7").Value = '  +26.23%  '

$ws.Range("D8").Value = '''0.399'
$ws.Range("E8").Value = '  +3.36%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").Value = '''0.893'
$ws.Range("E10").Value = '  +12.35%  '

$ws.Range("D11").Value = '3.397.79'
$ws.Range("E11").Value = '  +8.93%  '

$ws.Range("E12").Value = '  +1.66%  '

$ws.Range("D13").Value = '98.637.46'
$ws.Range("E13").Value = '  +2.30%  '

$ws.Range("D14").Value = '''36.51'
$ws.Range("E14").Value = '  +7.37%  '

$ws.Range("D15").Value = '''0.0000251'
$ws.Range("E15").Value = '  +4.48%  '

$ws.Range("D16").Value = '4.023.22'
$ws.Range("E16").Value = '  +8.62%  '

$ws.Range("E17").Value = '  +4.52%  '

$ws.Range("D18").Value = '3.385.38'
$ws.Range("E18").Value = '  +8.88%  '

$ws.Range("E19").Value = '  +1.55%  '

$ws.Range("D20").Value = '''15.34'
$ws.Range("E20").Value = '  +6.01%  '

$ws.Range("D21").Value = '''498.18'
$ws.Range("E21").Value = '  +1.67%  '

$ws.Range("D22").Value = '''6.26'
$ws.Range("E22").Value = '  +10.04%  '

$ws.Range("D23").Value = '''0.0000213'
$ws.Range("E23").Value = '  +9.75%  '

$ws.Range("D24").Value = '''9.46'
$ws.Range("E24").Value = '  +7.72%  '

$ws.Range("D25").Value = '''5.83'
$ws.Range("E25").Value = '  +5.49%  '

$ws.Range("D26").Value = '''89.53'
$ws.Range("E26").Value = '  +4.07%  '

$ws.Range("D27").Value = '''12.18'
$ws.Range("E27").Value = '  +4.47%  '

$ws.Range("E28").Value = '  +6.97%  '

$ws.Range("D29").Value = '''0.287'
$ws.Range("E29").Value = '  +20.72%  '

$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("E31").Value = '  +10.35%  '

$ws.Range("E32").Value = '  +7.29%  '

$ws.Range("D33").Value = '''9.72'
$ws.Range("E33").Value = '  +7.38%  '

$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +18.54%  '

$ws.Range("D35").Value = '''28.16'
$ws.Range("E35").Value = '  +7.20%  '

$ws.Range("D36").Value = '''7.43'
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").Value = '''2.00'
$ws.Range("E37").Value = '  +6.75%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.151'
$ws.Range("E38").Value = '  +0.36%  '

$ws.Range("D39").Value = '''0.476'
$ws.Range("E39").Value = '  +8.49%  '

$ws.Range("D40").Value = '''508.53'
$ws.Range("E40").Value = '  +3.54%  '

$ws.Range("E41").Value = '  +2.82%  '

$ws.Range("D42").Value = '''3.81'
$ws.Range("E42").Value = '  +5.12%  '

$ws.Range("E43").Value = '  +4.02%  '

$ws.Range("D44").Value = '''3.41'
$ws.Range("E44").Value = '  +6.19%  '

$ws.Range("D45").Value = '''0.792'
$ws.Range("E45").Value = '  +12.94%  '

$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").Value = '''160.30'
$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("E48").Value = '  +2.20%  '

$ws.Range("D49").Value = '''4.72'
$ws.Range("E49").Value = '  +8.64%  '

$ws.Range("D50").Value = '''46.74'
$ws.Range("E50").Value = '  +5.46%  '

$ws.Range("D51").Value = '''0.832'
$ws.Range("E51").Value = '  +13.38%  '

